$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date from 2021-05-06 to 2021-05-07
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2454491725761808
$ws.Range("E2").Value = 0.01116294568477172

$ws.Range("D3").Value = 0.5000848222095972
$ws.Range("E3").Value = 0.006537656903765621

$ws.Range("D4").Value = 0.09536442045329166
$ws.Range("E4").Value = 0.01130952380952399

$ws.Range("D5").Value = 0.102331442531942
$ws.Range("E5").Value = 0.009964987880419995

$ws.Range("D6").Value = 0.05677014222898843
$ws.Range("E6").Value = 0.01033475623455393

$ws.Range("E7").Value = 0.008694282121879393

# Restore sheet protection (sheet was protected before this edit)
$ws.Protect()
